$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# 1. The "_GoBack" bookmark (Word's "last edit location" marker) currently
#    sits at the end of the very first paragraph. Remove it from there.
# --------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --------------------------------------------------------------------------
# 2. The document's last paragraph is empty; it should end up containing the
#    text "222", split across two runs the way Word naturally would if "2"
#    was typed through an East-Asian IME (giving it an eastAsia font hint)
#    and "22" was then typed normally (no special run properties).
# --------------------------------------------------------------------------
$paraCountBefore = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($paraCountBefore)
$r = $lastPara.Range
$r.Collapse(1)

# A temporary sentinel character "X" is appended after the real text. This
# keeps the later bookmark-insertion point away from the paragraph-mark
# boundary (a position the host mishandles for collapsed bookmarks) - the
# sentinel is stripped again once the bookmark has been anchored.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t>22X</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# InsertXML inserts the new paragraph ahead of the original (still-empty)
# last paragraph instead of merging into it; fuse them back into one
# paragraph by deleting the paragraph mark between them.
$newPara = $d.Paragraphs.Item($paraCountBefore)
$mergeBoundary = $d.Range($newPara.Range.End - 1, $newPara.Range.End)
$mergeBoundary.Delete()

# --------------------------------------------------------------------------
# 3. Re-create "_GoBack" collapsed right after "222" (Word's usual spot for
#    it - the last place text was inserted), then remove the sentinel "X".
# --------------------------------------------------------------------------
$finalPara = $d.Paragraphs.Item($paraCountBefore)
$bmPos = $finalPara.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$finalPara2 = $d.Paragraphs.Item($paraCountBefore)
$xPos = $finalPara2.Range.End - 2
$sentinelRange = $d.Range($xPos, $xPos + 1)
$sentinelRange.Delete()
